# Atualiza os rótulos da linha de cabeçalho (linha 1) de cada planilha,
# prefixando os anos/intervalos com "Ano " ou "Intervalo " para que o
# Power BI reconheça automaticamente a primeira linha como cabeçalho.

$wb = $excel.ActiveWorkbook

# Planilhas 1, 2, 3, 5: colunas B1:E1 = anos (2015, 2030, 2040, 2050) -> "Ano <ano>"
$anoSheets = @(1, 2, 3, 5)
foreach ($idx in $anoSheets) {
    $ws = $wb.Worksheets.Item($idx)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $cell.Value2 = "Ano " + $cell.Value2
    }
}

# Planilha 4: colunas B1:E1 = períodos (2015, 2015-2030, 2031-2040, 2041-2050) -> "Intervalo <periodo>"
$ws4 = $wb.Worksheets.Item(4)
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws4.Range($col + "1")
    $cell.Value2 = "Intervalo " + $cell.Value2
}

# Planilha 6: apenas coluna B1 = 2015 -> "Ano 2015"
$ws6 = $wb.Worksheets.Item(6)
$cell6 = $ws6.Range("B1")
$cell6.Value2 = "Ano " + $cell6.Value2
